$wb = $excel.ActiveWorkbook

# Row -> new F-column value (applies to both "展览" and "全部类型" sheets)
$changes = @{
    2  = 1151
    3  = 866
    4  = 286
    8  = 2395
    9  = 7797
    12 = 390
    15 = 6
    17 = 8014
    19 = 1392
    24 = 332
    25 = 174
    30 = 429
    31 = 1162
    33 = 100
    38 = 71
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Cells.Item($row, 6).Value = $changes[$row]
    }
}
